$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the leads in the TEXT1 column (L2:L5) to the new test data.
$ws.Range("L2").Value = "dedic 1"
$ws.Range("L3").Value = "new nana 6"
$ws.Range("L4").Value = "adit saputra"
$ws.Range("L5").Value = "lead prospek kedua"

# Move the active selection from L6 to L5.
$ws.Range("L5").Select() | Out-Null
